# NEW_SC03 scripting: add rows 8 and 9 to the Talent_Acquisition sheet.
#
# Both rows mirror the existing row 7 layout (same business-unit, location,
# hire, payroll ... values) but represent two new scenarios:
#   row 8 -> CHANGE_SALARY_BASIS_FOR_EXISTING_EMP
#   row 9 -> EDIT_SALARY_PROPOSAL_REASON
# with their own scenario name, userName, personNumber and salaryAmount.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Talent_Acquisition")

# Pre-format the new rows as Text (numFmtId 49) up front, same as the rest of
# the sheet, so numeric-looking values (ids, amounts, zip codes, ...) are
# written as text rather than being auto-converted to numbers.
$ws.Range("B8:BB9").NumberFormat = "@"

# --- Row 8: CHANGE_SALARY_BASIS_FOR_EXISTING_EMP ---
$ws.Range("A8").Value = "CHANGE_SALARY_BASIS_FOR_EXISTING_EMP"
$ws.Range("AX8").Value = "3040301"
$ws.Range("AR8").Value = "3259840"
$ws.Range("AS8").Value = "80,000.00"
$ws.Range("B8").Value = "527 Johns Avenue"
$ws.Range("C8").Value = "Suite 436"
$ws.Range("D8").Value = "'1,000.00"
$ws.Range("E8").Value = "Full-time regular"
$ws.Range("F8").Value = "Finance"
$ws.Range("G8").Value = "Irving"
$ws.Range("H8").Value = "'510"
$ws.Range("I8").Value = "'61465"
$ws.Range("J8").Value = "United States"
$ws.Range("K8").Value = "United States 1"
$ws.Range("L8").Value = "Dallas"
$ws.Range("M8").Value = "12-Apr-1988"
$ws.Range("N8").Value = "Finance Internal Audit - EMEA"
$ws.Range("O8").Value = "Pearson Inc"
$ws.Range("P8").Value = "Bilingual Indicator"
$ws.Range("Q8").Value = "ADP Auto & Home"
$ws.Range("R8").Value = "uatnewhire67@abc.com"
$ws.Range("S8").Value = "Global Temporary Assignment"
$ws.Range("T8").Value = "GM Temporary Relocation"
$ws.Range("U8").Value = "13-Feb-2019"
$ws.Range("V8").Value = "'31"
$ws.Range("X8").Value = "Male"
$ws.Range("Y8").Value = "International Assignment"
$ws.Range("Z8").Value = "E"
$ws.Range("AA8").Value = "Hire"
$ws.Range("AB8").Value = "28-Jan-19"
$ws.Range("AC8").Value = "Additional Hire"
$ws.Range("AD8").Value = "Salaried"
$ws.Range("AE8").Value = "Director Audit & Compliance"
$ws.Range("AG8").Value = "NCS Pearson, Inc"
$ws.Range("AH8").Value = "AR-Buenos Aires-Humboldt 1509/13"
$ws.Range("AI8").Value = "Single"
$ws.Range("AK8").Value = "Social Security Number"
$ws.Range("AL8").Value = "12-Apr-2017"
$ws.Range("AM8").Value = "2717522"
$ws.Range("AN8").Value = "Welcome123"
$ws.Range("AO8").Value = "Yes"
$ws.Range("AP8").Value = "Bi-Weekly Exempt"
$ws.Range("AQ8").Value = "13-Feb-2019"
$ws.Range("AT8").Value = "US Annual Salary"
$ws.Range("AU8").Value = "TX"
$ws.Range("AV8").Value = "Work Phone"
$ws.Range("AW8").Value = "Home E-Mail"
$ws.Range("AY8").Value = "Not a Protected Veteran"
$ws.Range("AZ8").Value = "Employee"
$ws.Range("BA8").Value = "No"
$ws.Range("BB8").Value = "75038"

# --- Row 9: EDIT_SALARY_PROPOSAL_REASON ---
$ws.Range("A9").Value = "EDIT_SALARY_PROPOSAL_REASON"
$ws.Range("AX9").Value = "3259228"
$ws.Range("AR9").Value = "3259949"
$ws.Range("AS9").Value = "15.00"
$ws.Range("B9").Value = "527 Johns Avenue"
$ws.Range("C9").Value = "Suite 436"
$ws.Range("D9").Value = "'1,000.00"
$ws.Range("E9").Value = "Full-time regular"
$ws.Range("F9").Value = "Finance"
$ws.Range("G9").Value = "Irving"
$ws.Range("H9").Value = "'510"
$ws.Range("I9").Value = "'61465"
$ws.Range("J9").Value = "United States"
$ws.Range("K9").Value = "United States 1"
$ws.Range("L9").Value = "Dallas"
$ws.Range("M9").Value = "12-Apr-1988"
$ws.Range("N9").Value = "Finance Internal Audit - EMEA"
$ws.Range("O9").Value = "Pearson Inc"
$ws.Range("P9").Value = "Bilingual Indicator"
$ws.Range("Q9").Value = "ADP Auto & Home"
$ws.Range("R9").Value = "uatnewhire67@abc.com"
$ws.Range("S9").Value = "Global Temporary Assignment"
$ws.Range("T9").Value = "GM Temporary Relocation"
$ws.Range("U9").Value = "13-Feb-2019"
$ws.Range("V9").Value = "'31"
$ws.Range("X9").Value = "Male"
$ws.Range("Y9").Value = "International Assignment"
$ws.Range("Z9").Value = "E"
$ws.Range("AA9").Value = "Hire"
$ws.Range("AB9").Value = "28-Jan-19"
$ws.Range("AC9").Value = "Additional Hire"
$ws.Range("AD9").Value = "Salaried"
$ws.Range("AE9").Value = "Director Audit & Compliance"
$ws.Range("AG9").Value = "NCS Pearson, Inc"
$ws.Range("AH9").Value = "AR-Buenos Aires-Humboldt 1509/13"
$ws.Range("AI9").Value = "Single"
$ws.Range("AK9").Value = "Social Security Number"
$ws.Range("AL9").Value = "12-Apr-2017"
$ws.Range("AM9").Value = "2717522"
$ws.Range("AN9").Value = "Welcome123"
$ws.Range("AO9").Value = "Yes"
$ws.Range("AP9").Value = "Bi-Weekly Exempt"
$ws.Range("AQ9").Value = "13-Feb-2019"
$ws.Range("AT9").Value = "US Annual Salary"
$ws.Range("AU9").Value = "TX"
$ws.Range("AV9").Value = "Work Phone"
$ws.Range("AW9").Value = "Home E-Mail"
$ws.Range("AY9").Value = "Not a Protected Veteran"
$ws.Range("AZ9").Value = "Employee"
$ws.Range("BA9").Value = "No"
$ws.Range("BB9").Value = "75038"

# Match the saved selection / active cell from the authored workbook.
$ws.Range("AX9").Select()

